# Populate the "Steps" worksheet with the automation step table
# (STEP_ID / ACTION_NAME / RUN / DATASHEET / ITERATIONS / VALUE /
#  STORE_RESULT_AS / ON_FAILURE) describing the Open Browser / Close
# Browser steps, then apply the header highlight, column sizing and
# selection that come along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "STEP_ID"
$ws.Range("B1").Value = "ACTION_NAME"
$ws.Range("C1").Value = "RUN"
$ws.Range("D1").Value = "DATASHEET"
$ws.Range("E1").Value = "ITERATIONS"
$ws.Range("H1").Value = "ON_FAILURE"

# --- Step 10: Open Browser ------------------------------------------
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Open Browser"
$ws.Range("C2").Value = $true
$ws.Range("F1").Value = "VALUE"
$ws.Range("F2").Value = "Chrome"
$ws.Range("H2").Value = "ExitTest"
$ws.Range("D2").Value = "N/A"

# --- Step 20: Close Browser ------------------------------------------
$ws.Range("A3").Value = 20
$ws.Range("D3").Value = "N/A"
$ws.Range("H3").Value = "ExitTest"

$ws.Range("G1").Value = "STORE_RESULT_AS"
$ws.Range("B3").Value = "Close Browser"
$ws.Range("C3").Value = $true

# --- Formatting -------------------------------------------------------
# Highlight the new STORE_RESULT_AS header in maroon (FF800000)
$ws.Range("G1").Font.Color = 128

# Column widths to fit the new content
$ws.Columns.Item(2).ColumnWidth = 13.6    # ACTION_NAME
$ws.Columns.Item(4).ColumnWidth = 10.33   # DATASHEET
$ws.Columns.Item(5).ColumnWidth = 10.66   # ITERATIONS
$ws.Range("F1:G1").ColumnWidth = 10.66    # VALUE / STORE_RESULT_AS
$ws.Columns.Item(8).ColumnWidth = 11.13   # ON_FAILURE

# Leave the selection where the author left off
$ws.Range("H6").Select() | Out-Null
